# custom accuracy + 데이터 1000개
#
# Row 5 held one extra decimal digit of precision compared to the rest of
# the sheet's "custom accuracy" (numFmt 0.000 / 2-3 significant decimals).
# Normalize it to 2 decimal places, and drop the trailing data row (row 6)
# so the sheet ends at row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column -> rounded (2 decimal place) value for row 5.
$row5Values = [ordered]@{
    "B5"  = 9.5
    "C5"  = 6.98
    "D5"  = 0.87
    "E5"  = 20.9
    "F5"  = 16.65
    "G5"  = 7.42
    "H5"  = 29.84
    "I5"  = 11.63
    "J5"  = 5.06
    "K5"  = 7.36
    "L5"  = 8.37
    "M5"  = 8.96
    "N5"  = 2.42
    "O5"  = 7.52
    "P5"  = 10.6
    "Q5"  = 6.52
    "R5"  = 0.7
    "S5"  = 0.52
    "T5"  = 107.12
    "U5"  = 21.05
    "V5"  = 6.94
    "W5"  = 13.95
    "X5"  = 7.3
    "Y5"  = 1.28
    "Z5"  = 14.42
    "AA5" = 6.13
    "AB5" = 5.53
    "AC5" = 6.49
    "AD5" = 8.76
    "AE5" = 0.54
    "AF5" = 27.23
    "AG5" = 3.82
    "AH5" = 8.68
}

foreach ($addr in $row5Values.Keys) {
    $ws.Range($addr).Value = $row5Values[$addr]
}

# The sheet had 5 rows of readings plus a 6th that's being dropped as part
# of trimming the dataset; remove it entirely (dimension shrinks to AH5).
$ws.Rows.Item(6).Delete()
